$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") date updated from 2023-10-25 (45224) to 2023-11-03 (45233)
# for rows 2-5, keeping the existing cell formatting/style intact.
$newDate = Get-Date -Year 2023 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0

$ws.Range("C2:C5").Value = $newDate
